# Update "想去人数" (F column) figures for several events that were
# refreshed when the site data was regenerated.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets("展览")
$wsExpo.Range("F5").Value  = 858    # 吉安·COMIC LIFE次元假日06: 804 -> 858
$wsExpo.Range("F10").Value = 4708   # 南昌·萌卡动漫展: 4704 -> 4708
$wsExpo.Range("F20").Value = 3622   # 江西·JMG（江西广电）第二届UP动漫游戏博览会: 3618 -> 3622
$wsExpo.Range("F21").Value = 235    # 萍乡·AU10秋至国漫展: 203 -> 235
$wsExpo.Range("F32").Value = 791    # 南昌·CM04动漫游戏博览会: 787 -> 791
$wsExpo.Range("F33").Value = 2245   # 南昌·云芽动漫音乐嘉年华: 2239 -> 2245
$wsExpo.Range("F34").Value = 413    # 南昌·云芽动漫音乐嘉年华·封茗囧菌内场票: 411 -> 413

# Sheet "全部类型" (all types) - same events, rows offset by one due to an
# extra local row (31) not present in "展览".
$wsAll = $wb.Worksheets("全部类型")
$wsAll.Range("F5").Value  = 858     # 吉安·COMIC LIFE次元假日06: 805 -> 858
$wsAll.Range("F10").Value = 4708    # 南昌·萌卡动漫展: 4704 -> 4708
$wsAll.Range("F20").Value = 3622    # 江西·JMG（江西广电）第二届UP动漫游戏博览会: 3618 -> 3622
$wsAll.Range("F21").Value = 235     # 萍乡·AU10秋至国漫展: 203 -> 235
$wsAll.Range("F33").Value = 791     # 南昌·CM04动漫游戏博览会: 788 -> 791
$wsAll.Range("F34").Value = 2245    # 南昌·云芽动漫音乐嘉年华: 2239 -> 2245
$wsAll.Range("F35").Value = 413     # 南昌·云芽动漫音乐嘉年华·封茗囧菌内场票: 411 -> 413
